$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.638.13'
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").Value = '1.616.22'
$ws.Range("E3").Value = '  -0.82%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.989'
$ws.Range("E4").Value = '  -0.75%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.13'
$ws.Range("E5").Value = '  -1.16%  '

$ws.Range("E6").Value = '  -1.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.990'
$ws.Range("E7").Value = '  -0.66%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.03'
$ws.Range("E8").Value = '  -0.92%  '

$ws.Range("E9").Value = '  -1.59%  '

$ws.Range("E10").Value = '  -1.34%  '

$ws.Range("E11").Value = '  -0.75%  '

$ws.Range("D12").Value = '1.845.98'
$ws.Range("E12").Value = '  -0.72%  '

$ws.Range("D13").Value = '1.622.11'
$ws.Range("E13").Value = '  -0.41%  '

$ws.Range("E14").Value = '  -1.61%  '

$ws.Range("E15").Value = '  -1.46%  '

$ws.Range("E16").Value = '  -1.16%  '

$ws.Range("D17").Value = '27.654.48'
$ws.Range("E17").Value = '  -0.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '226.77'
$ws.Range("E18").Value = '  -1.60%  '

$ws.Range("E19").Value = '  +1.49%  '

$ws.Range("E20").Value = '  -1.02%  '

$ws.Range("E21").Value = '  -0.64%  '

$ws.Range("E22").Value = '  -1.33%  '

$ws.Range("E23").Value = '  -3.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.03'
$ws.Range("E24").Value = '  -1.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.35'
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("E26").Value = '  -1.07%  '

$ws.Range("E27").Value = '  -0.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.40'
$ws.Range("E28").Value = '  -1.23%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.990'
$ws.Range("E29").Value = '  -0.64%  '

$ws.Range("E30").Value = '  -1.13%  '

$ws.Range("E31").Value = '  -0.77%  '

$ws.Range("E32").Value = '  -1.44%  '

$ws.Range("E33").Value = '  -0.43%  '

$ws.Range("D34").Value = '1.388.50'
$ws.Range("E34").Value = '  -1.25%  '

$ws.Range("E35").Value = '  +1.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.994'
$ws.Range("E36").Value = '  -1.87%  '

$ws.Range("E37").Value = '  -1.44%  '

$ws.Range("E38").Value = '  +0.24%  '

$ws.Range("E39").Value = '  -0.77%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.839'
$ws.Range("E40").Value = '  -3.10%  '

$ws.Range("E41").Value = '  -1.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.990'
$ws.Range("E42").Value = '  -0.61%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.40'
$ws.Range("E43").Value = '  -1.77%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.81'
$ws.Range("E44").Value = '  -0.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.34'
$ws.Range("E45").Value = '  -2.98%  '

$ws.Range("D46").Value = '1.756.62'

$ws.Range("E47").Value = '  -3.64%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.50'
$ws.Range("E48").Value = '  -0.37%  '

$ws.Range("E49").Value = '  +1.01%  '

$ws.Range("E50").Value = '  -0.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.50'
$ws.Range("E51").Value = '  +0.57%  '
